# Add a new "2022-Q3" sheet before "2022-Q2", populate it with the
# new quarter's fund-holdings data, and update the "总计" summary sheet
# with a new leading row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before "2022-Q2"
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Header row (copy of the fund-detail sheet header)
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Helper data for the 11 fund rows (A column is a simple 0-based index)
$rows = @(
    @{ A=0;  B="012582"; C="交银施罗德品质增长一年持有期混合A";       D="49.48"; E="93.45"; F="5.82"; G="2.8797"; H=6 },
    @{ A=1;  B="010454"; C="交银施罗德内需增长一年持有期混合";         D="36.11"; E="93.41"; F="5.80"; G="2.0944"; H=7 },
    @{ A=2;  B="519714"; C="交银施罗德消费新驱动股票";                 D="34.98"; E="89.92"; F="5.71"; G="1.9974"; H=8 },
    @{ A=3;  B="004868"; C="交银施罗德股息优化混合";                   D="23.52"; E="86.53"; F="5.69"; G="1.3383"; H=7 },
    @{ A=4;  B="005004"; C="交银施罗德品质升级混合A";                  D="18.93"; E="90.92"; F="5.78"; G="1.0942"; H=8 },
    @{ A=5;  B="013882"; C="交银施罗德品质升级混合C";                  D="7.70";  E="90.92"; F="5.78"; G="0.4451"; H=8 },
    @{ A=6;  B="519710"; C="交银施罗德策略回报灵活配置混合";           D="6.77";  E="79.44"; F="6.10"; G="0.4130"; H=7 },
    @{ A=7;  B="012583"; C="交银施罗德品质增长一年持有期混合C";        D="2.37";  E="93.45"; F="5.82"; G="0.1379"; H=6 },
    @{ A=8;  B="161030"; C="富国中证体育产业指数A";                    D="1.59";  E="94.00"; F="5.17"; G="0.0822"; H=1 },
    @{ A=9;  B="013278"; C="富国中证体育产业指数C";                    D="0.42";  E="94.00"; F="5.17"; G="0.0217"; H=1 },
    @{ A=10; B="007943"; C="富安达中证 500 指数增强";                  D="0.23";  E="92.61"; F="1.22"; G="0.0028"; H=8 }
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row.A
    $q3.Cells.Item($r, 2).Value = "'" + $row.B
    $q3.Cells.Item($r, 3).Value = $row.C
    $q3.Cells.Item($r, 4).Value = "'" + $row.D
    $q3.Cells.Item($r, 5).Value = "'" + $row.E
    $q3.Cells.Item($r, 6).Value = "'" + $row.F
    $q3.Cells.Item($r, 7).Value = "'" + $row.G
    $q3.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new leading row for
#    2022-Q3 and shift the existing quarters down by one.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 11
$total.Cells.Item(2, 4).Value = 10.51

# Renumber the A-column index (0-based) for the rows that got shifted down
for ($row = 3; $row -le 9; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}
